$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that cannot be automated (set first so "Cannot be Automated" gets the
# lower shared-string index, matching the target workbook ordering)
$cannotAutomateRows = @(26, 27)
foreach ($r in $cannotAutomateRows) {
    $ws.Cells.Item($r, 10).Value = "Cannot be Automated"
}

# Rows that can be automated
$canAutomateRows = @(11, 12, 14, 15, 16, 22, 23, 25, 28, 29, 30)
foreach ($r in $canAutomateRows) {
    $ws.Cells.Item($r, 10).Value = "Can be Automated"
}

# The new, wider text ("Cannot be Automated" / "Can be Automated") needs the
# column to grow to fit, same as Excel's own best-fit recalculation.
$ws.Columns.Item(10).ColumnWidth = 18.33

# Re-select J30 to match final workbook selection state, and leave it as the active cell
$ws.Range("J30").Select() | Out-Null
